$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "*Project of your choice*user (views, model, controllers)*your topic*") {
        $r.Font.StrikeThrough = 1
    }
}
